$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This edit re-shuffles the (Id, Ost, Nord) triple -- columns A, Q, R --
# across rows 2-19, and additionally swaps the whole species record
# (columns B, D, E, F, G, H, K, L, M, N) between row 2 and row 5.
# ------------------------------------------------------------------

# Snapshot the "before" A/Q/R triples for every data row (2-19) first,
# since several of the writes below would otherwise clobber values we
# still need to read for a later row.
$idVal = @{}
$ostVal = @{}
$nordVal = @{}
for ($r = 2; $r -le 19; $r++) {
    $idVal[$r] = $ws.Cells.Item($r, 1).Value2
    $ostVal[$r] = $ws.Cells.Item($r, 17).Value2
    $nordVal[$r] = $ws.Cells.Item($r, 18).Value2
}

# Row -> row whose original (A,Q,R) triple it should receive.
$srcRow = @{
    2 = 5;  3 = 9;  4 = 19; 5 = 18; 6 = 11; 7 = 3;  8 = 4;  9 = 14; 10 = 12
    11 = 6; 12 = 7; 13 = 2; 14 = 10; 15 = 16; 16 = 17; 17 = 8; 18 = 13; 19 = 15
}

foreach ($r in $srcRow.Keys) {
    $src = $srcRow[$r]
    $ws.Cells.Item($r, 1).Value = $idVal[$src]
    $ws.Cells.Item($r, 17).Value = $ostVal[$src]
    $ws.Cells.Item($r, 18).Value = $nordVal[$src]
}

# ------------------------------------------------------------------
# Snapshot the species-record columns of rows 2 and 5 (B,D,E,F,G,H,K,L,M,N)
# before swapping them.
# ------------------------------------------------------------------
$cols = @(2, 4, 5, 6, 7, 8, 11, 12, 13, 14)  # B, D, E, F, G, H, K, L, M, N

$row2vals = @{}
$row5vals = @{}
foreach ($c in $cols) {
    $row2vals[$c] = $ws.Cells.Item(2, $c).Value2
    $row5vals[$c] = $ws.Cells.Item(5, $c).Value2
}

foreach ($c in $cols) {
    $v5 = $row5vals[$c]
    if ($null -eq $v5 -or $v5 -eq "") {
        $ws.Cells.Item(2, $c).Value = ""
    } else {
        $ws.Cells.Item(2, $c).Value = $v5
    }

    $v2 = $row2vals[$c]
    if ($null -eq $v2 -or $v2 -eq "") {
        $ws.Cells.Item(5, $c).Value = ""
    } else {
        $ws.Cells.Item(5, $c).Value = $v2
    }
}
